{"js": "// Replace every arithmetic-problem cell in the worksheet's single table\n// with its new expression. The table is 20 rows x 5 columns = 100 cells,\n// and every cell's text changes. We address cells positionally\n// (row, column) rather than via global text search/replace because a\n// few of the expression strings are substrings of one another (e.g.\n// \"5+19=\" is contained in \"35+19=\"), which would make a blind\n// find/replace ambiguous. Using `cell.value = \u2026` (rather than\n// `body.insertText(\u2026, \"Replace\")`) keeps each cell's existing run\n// formatting (TimeNewRoman, size 30) untouched \u2014 only the `<w:t>`\n// content changes, matching the source edit.\nconst newGrid = [\n  [\"99-39=\", \"5+3=\", \"28+33=\", \"55-29=\", \"93-19=\"],\n  [\"47-15=\", \"80-50=\", \"78-32=\", \"87-42=\", \"8+21=\"],\n  [\"63-29=\", \"88-85=\", \"3+63=\", \"19-11=\", \"64+33=\"],\n  [\"38+56=\", \"27-5=\", \"55+20=\", \"62-22=\", \"16-14=\"],\n  [\"36-18=\", \"60-6=\", \"50+12=\", \"49-16=\", \"6+60=\"],\n  [\"43-15=\", \"72-25=\", \"50+41=\", \"52+15=\", \"8+69=\"],\n  [\"17+81=\", \"47-34=\", \"40+12=\", \"75+24=\", \"57-50=\"],\n  [\"69-19=\", \"32+0=\", \"47+39=\", \"68-36=\", \"29+2=\"],\n  [\"61-41=\", \"55-12=\", \"20+19=\", \"4+89=\", \"65+18=\"],\n  [\"79-18=\", \"25-12=\", \"69-63=\", \"26+21=\", \"64+4=\"],\n  [\"11+42=\", \"72-4=\", \"43+39=\", \"4+36=\", \"32+1=\"],\n  [\"15+2=\", \"69+3=\", \"19+42=\", \"26+55=\", \"48-30=\"],\n  [\"80-13=\", \"63-25=\", \"45-10=\", \"87-82=\", \"63+19=\"],\n  [\"46-30=\", \"35-19=\", \"53-45=\", \"81-22=\", \"9+20=\"],\n  [\"24+36=\", \"15-5=\", \"37+0=\", \"93-26=\", \"72-50=\"],\n  [\"84-80=\", \"92-88=\", \"95-74=\", \"74-20=\", \"68-32=\"],\n  [\"82-25=\", \"96-9=\", \"3+1=\", \"3+82=\", \"23+29=\"],\n  [\"91-64=\", \"31+10=\", \"3+71=\", \"81+13=\", \"83-18=\"],\n  [\"41-40=\", \"25-11=\", \"2+68=\", \"20-9=\", \"15+49=\"],\n  [\"7+28=\", \"96-18=\", \"22+17=\", \"6+71=\", \"45-8=\"]\n];\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nif (tables.items.length === 0) {\n  throw new Error(\"Expected a table in the document body, found none.\");\n}\n\nconst table = tables.items[0];\n\nfor (let r = 0; r < newGrid.length; r++) {\n  const row = newGrid[r];\n  for (let c = 0; c < row.length; c++) {\n    const cell = table.getCell(r, c);\n    cell.value = row[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Replace every arithmetic-problem cell in the worksheet's single table\n# with its new expression. The table is 20 rows x 5 columns = 100 cells,\n# and every cell's text changes. Cells are addressed positionally\n# (row, column) rather than via a global Find/Replace because a few of\n# the expression strings are substrings of one another (e.g. \"5+19=\"\n# is contained in \"35+19=\"), which would make a blind find/replace\n# ambiguous. Assigning directly to `Cell.Range.Text` (rather than using\n# Find.Execute's Replacement) keeps each cell's existing run formatting\n# (TimeNewRoman, size 30) untouched \u2014 only the text content changes,\n# matching the source edit.\n\n$d = $word.ActiveDocument\n\n$newGrid = @(\n  @(\"99-39=\",\"5+3=\",\"28+33=\",\"55-29=\",\"93-19=\"),\n  @(\"47-15=\",\"80-50=\",\"78-32=\",\"87-42=\",\"8+21=\"),\n  @(\"63-29=\",\"88-85=\",\"3+63=\",\"19-11=\",\"64+33=\"),\n  @(\"38+56=\",\"27-5=\",\"55+20=\",\"62-22=\",\"16-14=\"),\n  @(\"36-18=\",\"60-6=\",\"50+12=\",\"49-16=\",\"6+60=\"),\n  @(\"43-15=\",\"72-25=\",\"50+41=\",\"52+15=\",\"8+69=\"),\n  @(\"17+81=\",\"47-34=\",\"40+12=\",\"75+24=\",\"57-50=\"),\n  @(\"69-19=\",\"32+0=\",\"47+39=\",\"68-36=\",\"29+2=\"),\n  @(\"61-41=\",\"55-12=\",\"20+19=\",\"4+89=\",\"65+18=\"),\n  @(\"79-18=\",\"25-12=\",\"69-63=\",\"26+21=\",\"64+4=\"),\n  @(\"11+42=\",\"72-4=\",\"43+39=\",\"4+36=\",\"32+1=\"),\n  @(\"15+2=\",\"69+3=\",\"19+42=\",\"26+55=\",\"48-30=\"),\n  @(\"80-13=\",\"63-25=\",\"45-10=\",\"87-82=\",\"63+19=\"),\n  @(\"46-30=\",\"35-19=\",\"53-45=\",\"81-22=\",\"9+20=\"),\n  @(\"24+36=\",\"15-5=\",\"37+0=\",\"93-26=\",\"72-50=\"),\n  @(\"84-80=\",\"92-88=\",\"95-74=\",\"74-20=\",\"68-32=\"),\n  @(\"82-25=\",\"96-9=\",\"3+1=\",\"3+82=\",\"23+29=\"),\n  @(\"91-64=\",\"31+10=\",\"3+71=\",\"81+13=\",\"83-18=\"),\n  @(\"41-40=\",\"25-11=\",\"2+68=\",\"20-9=\",\"15+49=\"),\n  @(\"7+28=\",\"96-18=\",\"22+17=\",\"6+71=\",\"45-8=\")\n)\n\n$table = $d.Tables.Item(1)\n\nfor ($r = 0; $r -lt $newGrid.Length; $r++) {\n  $row = $newGrid[$r]\n  for ($c = 0; $c -lt $row.Length; $c++) {\n    $cell = $table.Cell($r + 1, $c + 1)\n    $cell.Range.Text = $row[$c]\n  }\n}\n"}
